$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: keep Output7 in A, move the old "name" header (Output1) to B,
# and declare a new Output2 header in C for the index-check results.
$ws.Range("A1").Value = "Output7"
$ws.Range("B1").Value = "Output1"
$ws.Range("C1").Value = "Output2"

# Row 2
$ws.Range("A2").Value = "John"
$ws.Range("B2").Value = "Josh"
$ws.Range("C2").Value = 5

# Row 3
$ws.Range("A3").Value = "Viridian"
$ws.Range("B3").Value = "Merti"
$ws.Range("C3").Value = 6

# Row 4
$ws.Range("A4").Value = "Merti"
$ws.Range("B4").Value = "Viridian"
$ws.Range("C4").Value = 43

# Row 5
$ws.Range("A5").Value = "Pepe"
$ws.Range("B5").Value = "John"
$ws.Range("C5").Value = 48

# Row 6
$ws.Range("A6").Value = "Nadia"
$ws.Range("B6").Value = "Nadia"
$ws.Range("C6").Value = 6

# Row 7
$ws.Range("A7").Value = "Josh"
$ws.Range("B7").Value = "Pepe"
$ws.Range("C7").Value = 8
